$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.247.26'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '3.896.30'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.612'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.53%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.716'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.35%  '
$ws.Range('E10').Value = '  -5.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000336'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('D13').Value = '4.522.22'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '3.915.91'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.135'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').Value = '69.210.20'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '425.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('E22').Value = '  -6.59%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '88.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.72%  '
$ws.Range('E25').Value = '  +10.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '680.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.00%  '
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('E31').Value = '  -3.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '68.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.85%  '
$ws.Range('D34').Value = '0.0₃0870'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('E35').Value = '  +8.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '40.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.149'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.30%  '
$ws.Range('E42').Value = '  -3.37%  '
$ws.Range('E43').Value = '  +7.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('E47').Value = '  +18.12%  '
$ws.Range('E48').Value = '  +6.59%  '
$ws.Range('D49').Value = '0.0₆0351'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.743.50'
$ws.Range('E51').Value = '  +14.20%  '
